$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (R and S) right after the current last data column (Q),
# cloning the column formatting (this mirrors a user extending the yearly table
# with two more year columns, 2021 and 2022).
$ws.Range("R1:S12").Insert(-4161)

# --- Row 4 (year headers) ---
$ws.Range("R4").Value = 2021
$ws.Range("S4").Value = 2022

# --- Row 5 ---
$ws.Range("P5").Value = 25.6
$ws.Range("Q5").Value = 23.8
$ws.Range("R5").Value = 26.8
$ws.Range("S5").Value = 26.8

# --- Row 6 ---
$ws.Range("P6").Value = 18.600000000000001
$ws.Range("Q6").Value = 16.7
$ws.Range("R6").Value = 19.3
$ws.Range("S6").Value = 19.3

# --- Row 7 ---
$ws.Range("R7").Value = "-"
$ws.Range("S7").Value = "-"

# --- Row 8 ---
$ws.Range("P8").Value = 2.1
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 1.8
$ws.Range("S8").Value = 1.8

# --- Row 9 ---
$ws.Range("P9").Value = 4.9000000000000004
$ws.Range("Q9").Value = 5.2
$ws.Range("R9").Value = 5.7
$ws.Range("S9").Value = 5.7

# --- Row 10 ---
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0

# Update the active selection as left after the edit.
[void]$ws.Range("T3").Select()
